$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update each row (2-21) with the re-sorted team data plus new Rodada 31 (AF) / Rodada 32 (AG) scores.
# Row 2: SHOWZES F C
$ws.Range("A2").Value = "SHOWZES F C"
$row2 = New-Object 'object[,]' 1,32
$row2[0,0] = 73.260000000000005
$row2[0,1] = 92.79
$row2[0,2] = 71.02
$row2[0,3] = 99.12
$row2[0,4] = 72.88
$row2[0,5] = 105.36
$row2[0,6] = 117.59
$row2[0,7] = 110
$row2[0,8] = 79.92
$row2[0,9] = 65.069999999999993
$row2[0,10] = 126.77
$row2[0,11] = 68.36
$row2[0,12] = 117.97
$row2[0,13] = 66.8
$row2[0,14] = 107.8
$row2[0,15] = 70.67
$row2[0,16] = 56.81
$row2[0,17] = 93.47
$row2[0,18] = 62.55
$row2[0,19] = 75.38
$row2[0,20] = 193.46
$row2[0,21] = 77.56
$row2[0,22] = 98
$row2[0,23] = 71.28
$row2[0,24] = 53.85
$row2[0,25] = 92.44
$row2[0,26] = 99.07
$row2[0,27] = 116.54
$row2[0,28] = 108.58
$row2[0,29] = 35.39
$row2[0,30] = 90.47
$row2[0,31] = 124.84
$ws.Range("B2:AG2").Value = $row2

# Row 3: GD99.FC
$ws.Range("A3").Value = "GD99.FC "
$row3 = New-Object 'object[,]' 1,32
$row3[0,0] = 62.53
$row3[0,1] = 89.88
$row3[0,2] = 74.52
$row3[0,3] = 121.4
$row3[0,4] = 68.13
$row3[0,5] = 76.11
$row3[0,6] = 95.55
$row3[0,7] = 133.03
$row3[0,8] = 61.62
$row3[0,9] = 73.17
$row3[0,10] = 92.47
$row3[0,11] = 75.17
$row3[0,12] = 94.27
$row3[0,13] = 71.45
$row3[0,14] = 115.1
$row3[0,15] = 76.599999999999994
$row3[0,16] = 66.75
$row3[0,17] = 94.77
$row3[0,18] = 68.599999999999994
$row3[0,19] = 82.16
$row3[0,20] = 199.71
$row3[0,21] = 91.96
$row3[0,22] = 106.5
$row3[0,23] = 85.18
$row3[0,24] = 79.349999999999994
$row3[0,25] = 93.94
$row3[0,26] = 91.47
$row3[0,27] = 92.84
$row3[0,28] = 108.53
$row3[0,29] = 23.99
$row3[0,30] = 96.17
$row3[0,31] = 98.94
$ws.Range("B3:AG3").Value = $row3

# Row 4: SANTOS M.T.T
$ws.Range("A4").Value = "SANTOS M.T.T"
$row4 = New-Object 'object[,]' 1,32
$row4[0,0] = 63.85
$row4[0,1] = 72.69
$row4[0,2] = 71.62
$row4[0,3] = 93.68
$row4[0,4] = 75.62
$row4[0,5] = 110.21
$row4[0,6] = 91.99
$row4[0,7] = 152.5
$row4[0,8] = 94.1
$row4[0,9] = 70.569999999999993
$row4[0,10] = 102.97
$row4[0,11] = 81.96
$row4[0,12] = 114.91
$row4[0,13] = 77.900000000000006
$row4[0,14] = 100
$row4[0,15] = 63.57
$row4[0,16] = 42.61
$row4[0,17] = 93.1
$row4[0,18] = 89.05
$row4[0,19] = 86.06
$row4[0,20] = 180.01
$row4[0,21] = 65.09
$row4[0,22] = 110.47
$row4[0,23] = 45.83
$row4[0,24] = 51.75
$row4[0,25] = 55.55
$row4[0,26] = 105.97
$row4[0,27] = 127.14
$row4[0,28] = 112.53
$row4[0,29] = 64.19
$row4[0,30] = 88.52
$row4[0,31] = 103.54
$ws.Range("B4:AG4").Value = $row4

# Row 5: CRV ROCHA
$ws.Range("A5").Value = "CRV ROCHA "
$row5 = New-Object 'object[,]' 1,32
$row5[0,0] = 74.25
$row5[0,1] = 78.290000000000006
$row5[0,2] = 69.3
$row5[0,3] = 97.46
$row5[0,4] = 72.22
$row5[0,5] = 104.86
$row5[0,6] = 89.29
$row5[0,7] = 122.35
$row5[0,8] = 77.62
$row5[0,9] = 63.12
$row5[0,10] = 106.01
$row5[0,11] = 77.16
$row5[0,12] = 119.31
$row5[0,13] = 55.5
$row5[0,14] = 96.3
$row5[0,15] = 74.08
$row5[0,16] = 62.81
$row5[0,17] = 108.9
$row5[0,18] = 78.099999999999994
$row5[0,19] = 74.61
$row5[0,20] = 187.51
$row5[0,21] = 86.56
$row5[0,22] = 66.45
$row5[0,23] = 69.45
$row5[0,24] = 54.16
$row5[0,25] = 60.62
$row5[0,26] = 114.87
$row5[0,27] = 133.63999999999999
$row5[0,28] = 122.58
$row5[0,29] = 55.79
$row5[0,30] = 84.12
$row5[0,31] = 78.290000000000006
$ws.Range("B5:AG5").Value = $row5

# Row 6: MARXU
$ws.Range("A6").Value = "MARXU"
$row6 = New-Object 'object[,]' 1,32
$row6[0,0] = 90.29
$row6[0,1] = 87.29
$row6[0,2] = 85.05
$row6[0,3] = 115.51
$row6[0,4] = 56.68
$row6[0,5] = 106.56
$row6[0,6] = 75.91
$row6[0,7] = 116.94
$row6[0,8] = 90.67
$row6[0,9] = 66.62
$row6[0,10] = 118.14
$row6[0,11] = 44.66
$row6[0,12] = 97.66
$row6[0,13] = 65.400000000000006
$row6[0,14] = 112.3
$row6[0,15] = 71.86
$row6[0,16] = 51.21
$row6[0,17] = 86.66
$row6[0,18] = 62.3
$row6[0,19] = 72.86
$row6[0,20] = 191.06
$row6[0,21] = 83.06
$row6[0,22] = 69.95
$row6[0,23] = 49.81
$row6[0,24] = 68.150000000000006
$row6[0,25] = 76.94
$row6[0,26] = 64.97
$row6[0,27] = 134.94
$row6[0,28] = 74.790000000000006
$row6[0,29] = 37.590000000000003
$row6[0,30] = 84.52
$row6[0,31] = 106.14
$ws.Range("B6:AG6").Value = $row6

# Row 7: ANAEL O CRUEL
$ws.Range("A7").Value = "ANAEL O CRUEL"
$row7 = New-Object 'object[,]' 1,32
$row7[0,0] = 74.489999999999995
$row7[0,1] = 55.42
$row7[0,2] = 73.17
$row7[0,3] = 91.51
$row7[0,4] = 90.61
$row7[0,5] = 78.430000000000007
$row7[0,6] = 76.44
$row7[0,7] = 108.89
$row7[0,8] = 84.37
$row7[0,9] = 76.77
$row7[0,10] = 102.35
$row7[0,11] = 61.14
$row7[0,12] = 122.04
$row7[0,13] = 68.7
$row7[0,14] = 107.1
$row7[0,15] = 59.7
$row7[0,16] = 41.16
$row7[0,17] = 114.9
$row7[0,18] = 64.599999999999994
$row7[0,19] = 68.56
$row7[0,20] = 168.76
$row7[0,21] = 77.959999999999994
$row7[0,22] = 99.2
$row7[0,23] = 39.58
$row7[0,24] = 59.86
$row7[0,25] = 86.04
$row7[0,26] = 49.97
$row7[0,27] = 133.13999999999999
$row7[0,28] = 92.88
$row7[0,29] = 54.79
$row7[0,30] = 101.87
$row7[0,31] = 114.54
$ws.Range("B7:AG7").Value = $row7

# Row 8: PRESTINI FC
$ws.Range("A8").Value = "PRESTINI FC"
$row8 = New-Object 'object[,]' 1,32
$row8[0,0] = 74.790000000000006
$row8[0,1] = 82.29
$row8[0,2] = 56.1
$row8[0,3] = 118.28
$row8[0,4] = 88.43
$row8[0,5] = 74.06
$row8[0,6] = 87.86
$row8[0,7] = 110.44
$row8[0,8] = 64.42
$row8[0,9] = 59.77
$row8[0,10] = 96.72
$row8[0,11] = 56.46
$row8[0,12] = 117.25
$row8[0,13] = 60.8
$row8[0,14] = 99
$row8[0,15] = 81.66
$row8[0,16] = 69.41
$row8[0,17] = 94
$row8[0,18] = 72.150000000000006
$row8[0,19] = 58.36
$row8[0,20] = 188.61
$row8[0,21] = 83.76
$row8[0,22] = 100.47
$row8[0,23] = 43.63
$row8[0,24] = 54.2
$row8[0,25] = 98.84
$row8[0,26] = 96.45
$row8[0,27] = 131.04
$row8[0,28] = 109.43
$row8[0,29] = 57.3
$row8[0,30] = 76.97
$row8[0,31] = 11.74
$ws.Range("B8:AG8").Value = $row8

# Row 9: ALTO BRACINHO F.C
$ws.Range("A9").Value = "ALTO BRACINHO F.C "
$row9 = New-Object 'object[,]' 1,32
$row9[0,0] = 70.55
$row9[0,1] = 106.22
$row9[0,2] = 58.2
$row9[0,3] = 101.96
$row9[0,4] = 91.52
$row9[0,5] = 97.81
$row9[0,6] = 104.24
$row9[0,7] = 145.44999999999999
$row9[0,8] = 66.52
$row9[0,9] = 49.1
$row9[0,10] = 81.87
$row9[0,11] = 63.41
$row9[0,12] = 107.41
$row9[0,13] = 63.65
$row9[0,14] = 91.8
$row9[0,15] = 73.28
$row9[0,16] = 65.31
$row9[0,17] = 68.67
$row9[0,18] = 75.45
$row9[0,19] = 75.86
$row9[0,20] = 172.91
$row9[0,21] = 86.11
$row9[0,22] = 91.05
$row9[0,23] = 71.58
$row9[0,24] = 40.6
$row9[0,25] = 91.34
$row9[0,26] = 75.87
$row9[0,27] = 90.96
$row9[0,28] = 70.180000000000007
$row9[0,29] = 56
$row9[0,30] = 77.47
$row9[0,31] = 81.040000000000006
$ws.Range("B9:AG9").Value = $row9

# Row 10: ISAR7
$ws.Range("A10").Value = "ISAR7"
$row10 = New-Object 'object[,]' 1,32
$row10[0,0] = 70.25
$row10[0,1] = 66.69
$row10[0,2] = 64.819999999999993
$row10[0,3] = 98.8
$row10[0,4] = 77.53
$row10[0,5] = 84.25
$row10[0,6] = 89.49
$row10[0,7] = 92.21
$row10[0,8] = 86.53
$row10[0,9] = 63.95
$row10[0,10] = 99.46
$row10[0,11] = 52.56
$row10[0,12] = 103.11
$row10[0,13] = 53.7
$row10[0,14] = 90.49
$row10[0,15] = 76.58
$row10[0,16] = 55.85
$row10[0,17] = 96.68
$row10[0,18] = 66.31
$row10[0,19] = 68.650000000000006
$row10[0,20] = 154.22
$row10[0,21] = 59.89
$row10[0,22] = 67.23
$row10[0,23] = 58.61
$row10[0,24] = 58.75
$row10[0,25] = 72.62
$row10[0,26] = 101.99
$row10[0,27] = 107.06
$row10[0,28] = 88.43
$row10[0,29] = 56.04
$row10[0,30] = 64.27
$row10[0,31] = 50.22
$ws.Range("B10:AG10").Value = $row10

# Row 11: WILIAN H. FC
$ws.Range("A11").Value = "WILIAN H. FC"
$row11 = New-Object 'object[,]' 1,32
$row11[0,0] = 68.94
$row11[0,1] = 43.43
$row11[0,2] = 46.82
$row11[0,3] = 59.22
$row11[0,4] = 66.08
$row11[0,5] = 102.56
$row11[0,6] = 91.49
$row11[0,7] = 85.54
$row11[0,8] = 62.98
$row11[0,9] = 85.85
$row11[0,10] = 81.69
$row11[0,11] = 62.16
$row11[0,12] = 91.47
$row11[0,13] = 56.45
$row11[0,14] = 42.52
$row11[0,15] = 55.78
$row11[0,16] = 53.11
$row11[0,17] = 100.67
$row11[0,18] = 79.61
$row11[0,19] = 91.61
$row11[0,20] = 169.87
$row11[0,21] = 88.54
$row11[0,22] = 91.8
$row11[0,23] = 62.84
$row11[0,24] = 58.5
$row11[0,25] = 78.13
$row11[0,26] = 68.92
$row11[0,27] = 132.86000000000001
$row11[0,28] = 122.6
$row11[0,29] = 27.69
$row11[0,30] = 70.02
$row11[0,31] = 86.55
$ws.Range("B11:AG11").Value = $row11

# Row 12: TAKAMASSA NOMUROO
$ws.Range("A12").Value = "TAKAMASSA NOMUROO"
$row12 = New-Object 'object[,]' 1,32
$row12[0,0] = 56.11
$row12[0,1] = 91.42
$row12[0,2] = 82.12
$row12[0,3] = 70.73
$row12[0,4] = 39.99
$row12[0,5] = 104.45
$row12[0,6] = 105.44
$row12[0,7] = 52.03
$row12[0,8] = 45.67
$row12[0,9] = 67.87
$row12[0,10] = 65.69
$row12[0,11] = 10.87
$row12[0,12] = 94.21
$row12[0,13] = 57.4
$row12[0,14] = 97.2
$row12[0,15] = 65.53
$row12[0,16] = 84.91
$row12[0,17] = 65.099999999999994
$row12[0,18] = 63.5
$row12[0,19] = 76.86
$row12[0,20] = 172.91
$row12[0,21] = 56.26
$row12[0,22] = 65.02
$row12[0,23] = 80.86
$row12[0,24] = 51.3
$row12[0,25] = 89.94
$row12[0,26] = 86.92
$row12[0,27] = 109.44
$row12[0,28] = 91.43
$row12[0,29] = 37.44
$row12[0,30] = 84.97
$row12[0,31] = 105.04
$ws.Range("B12:AG12").Value = $row12

# Row 13: E. C. HOCHLEITNER
$ws.Range("A13").Value = "E. C. HOCHLEITNER"
$row13 = New-Object 'object[,]' 1,32
$row13[0,0] = 63.55
$row13[0,1] = 79.790000000000006
$row13[0,2] = 71.900000000000006
$row13[0,3] = 82.81
$row13[0,4] = 68.92
$row13[0,5] = 92.16
$row13[0,6] = 81.540000000000006
$row13[0,7] = 95.55
$row13[0,8] = 31.22
$row13[0,9] = 36.700000000000003
$row13[0,10] = 115.42
$row13[0,11] = 7.8
$row13[0,12] = 96.71
$row13[0,13] = 73.5
$row13[0,14] = 69.180000000000007
$row13[0,15] = 81.88
$row13[0,16] = 62.92
$row13[0,17] = 66.17
$row13[0,18] = 62.8
$row13[0,19] = 65.55
$row13[0,20] = 164.56
$row13[0,21] = 80.16
$row13[0,22] = 59.45
$row13[0,23] = 43.35
$row13[0,24] = 51.65
$row13[0,25] = 91.34
$row13[0,26] = 61.38
$row13[0,27] = 111.94
$row13[0,28] = 115.69
$row13[0,29] = 43.14
$row13[0,30] = 78.569999999999993
$row13[0,31] = 69.77
$ws.Range("B13:AG13").Value = $row13

# Row 14: JARAGUÁ 99 EC
$ws.Range("A14").Value = "JARAGUÁ 99 EC"
$row14 = New-Object 'object[,]' 1,32
$row14[0,0] = 0
$row14[0,1] = 58.59
$row14[0,2] = 88.52
$row14[0,3] = 114.61
$row14[0,4] = 80.83
$row14[0,5] = 82.56
$row14[0,6] = 89.95
$row14[0,7] = 106
$row14[0,8] = 54.27
$row14[0,9] = 90.64
$row14[0,10] = 89.77
$row14[0,11] = 61.26
$row14[0,12] = 80.97
$row14[0,13] = 53.7
$row14[0,14] = 75.819999999999993
$row14[0,15] = 53.65
$row14[0,16] = 64.25
$row14[0,17] = 82.42
$row14[0,18] = 78.400000000000006
$row14[0,19] = 68.36
$row14[0,20] = 124.21
$row14[0,21] = 95.61
$row14[0,22] = 91.1
$row14[0,23] = 39.06
$row14[0,24] = 58.25
$row14[0,25] = 80.84
$row14[0,26] = 67.08
$row14[0,27] = 95.44
$row14[0,28] = 41.59
$row14[0,29] = 41.14
$row14[0,30] = 57.79
$row14[0,31] = 81.64
$ws.Range("B14:AG14").Value = $row14

# Row 15: BARZEA51
$ws.Range("A15").Value = "BARZEA51 "
$row15 = New-Object 'object[,]' 1,32
$row15[0,0] = 89.63
$row15[0,1] = 102.92
$row15[0,2] = 73.319999999999993
$row15[0,3] = 132.21
$row15[0,4] = 85.18
$row15[0,5] = 82.41
$row15[0,6] = 89.74
$row15[0,7] = 95.19
$row15[0,8] = 80.63
$row15[0,9] = 29.96
$row15[0,10] = 97.47
$row15[0,11] = 65.36
$row15[0,12] = 86.21
$row15[0,13] = 16.5
$row15[0,14] = 105.5
$row15[0,15] = 51.67
$row15[0,16] = 59.74
$row15[0,17] = 69.8
$row15[0,18] = 48.31
$row15[0,19] = 39.81
$row15[0,20] = 161.47
$row15[0,21] = 73.489999999999995
$row15[0,22] = 78.900000000000006
$row15[0,23] = 63.31
$row15[0,24] = 46.15
$row15[0,25] = 50.9
$row15[0,26] = 25.45
$row15[0,27] = 85.18
$row15[0,28] = 93.04
$row15[0,29] = 31.93
$row15[0,30] = 66.55
$row15[0,31] = 63.69
$ws.Range("B15:AG15").Value = $row15

# Row 16: GRÊMIO FOOT-BALL JARAGUAENSE
$ws.Range("A16").Value = "GRÊMIO FOOT-BALL JARAGUAENSE"
$row16 = New-Object 'object[,]' 1,32
$row16[0,0] = 84.38
$row16[0,1] = 72.319999999999993
$row16[0,2] = 63.92
$row16[0,3] = 87.62
$row16[0,4] = 89.88
$row16[0,5] = 112.64
$row16[0,6] = 90.2
$row16[0,7] = 73.959999999999994
$row16[0,8] = 68.03
$row16[0,9] = 82.72
$row16[0,10] = 125.92
$row16[0,11] = 72.5
$row16[0,12] = 55.8
$row16[0,13] = 69.849999999999994
$row16[0,14] = 67.3
$row16[0,15] = 68.7
$row16[0,16] = 51.15
$row16[0,17] = 81.36
$row16[0,18] = 44.58
$row16[0,19] = 78.13
$row16[0,20] = 150.11000000000001
$row16[0,21] = 84.41
$row16[0,22] = 93
$row16[0,23] = 19.41
$row16[0,24] = 31.6
$row16[0,25] = 52.82
$row16[0,26] = 43.07
$row16[0,27] = 92.55
$row16[0,28] = 77.06
$row16[0,29] = 28.69
$row16[0,30] = 20.02
$row16[0,31] = 88.92
$ws.Range("B16:AG16").Value = $row16

# Row 17: VICTEAM
$ws.Range("A17").Value = "VICTEAM"
$row17 = New-Object 'object[,]' 1,32
$row17[0,0] = 0
$row17[0,1] = 48.37
$row17[0,2] = 50.22
$row17[0,3] = 63.98
$row17[0,4] = 60.78
$row17[0,5] = 80.78
$row17[0,6] = 84.36
$row17[0,7] = 113.25
$row17[0,8] = 34.950000000000003
$row17[0,9] = 80.12
$row17[0,10] = 69.31
$row17[0,11] = 51.66
$row17[0,12] = 94.3
$row17[0,13] = 59.9
$row17[0,14] = 90.3
$row17[0,15] = 61.17
$row17[0,16] = 69.209999999999994
$row17[0,17] = 93.97
$row17[0,18] = 60.75
$row17[0,19] = 44.15
$row17[0,20] = 162.46
$row17[0,21] = 77.260000000000005
$row17[0,22] = 59.07
$row17[0,23] = 49.08
$row17[0,24] = 53.14
$row17[0,25] = 37.19
$row17[0,26] = 30.75
$row17[0,27] = 49.98
$row17[0,28] = 60.15
$row17[0,29] = 19.25
$row17[0,30] = 49.54
$row17[0,31] = 53.41
$ws.Range("B17:AG17").Value = $row17

# Row 18: SIMPLICIO SPORT CLUB
$ws.Range("A18").Value = "SIMPLICIO SPORT CLUB"
$row18 = New-Object 'object[,]' 1,32
$row18[0,0] = 55.38
$row18[0,1] = 55.74
$row18[0,2] = 47.97
$row18[0,3] = 84.66
$row18[0,4] = 96.86
$row18[0,5] = 79.5
$row18[0,6] = 60.03
$row18[0,7] = 104.98
$row18[0,8] = 67.349999999999994
$row18[0,9] = 68.31
$row18[0,10] = 60.42
$row18[0,11] = 73.709999999999994
$row18[0,12] = 55.09
$row18[0,13] = 37.65
$row18[0,14] = 62.42
$row18[0,15] = 51.5
$row18[0,16] = 46.3
$row18[0,17] = 75.900000000000006
$row18[0,18] = 87.36
$row18[0,19] = 82.28
$row18[0,20] = 160.36000000000001
$row18[0,21] = 50.36
$row18[0,22] = 61.6
$row18[0,23] = 46.11
$row18[0,24] = 48.45
$row18[0,25] = 39.07
$row18[0,26] = 27.02
$row18[0,27] = 54.65
$row18[0,28] = 47.71
$row18[0,29] = 49.69
$row18[0,30] = 31.62
$row18[0,31] = 27.68
$ws.Range("B18:AG18").Value = $row18

# Row 19: MORUMBIHEXA
$ws.Range("A19").Value = "MORUMBIHEXA"
$row19 = New-Object 'object[,]' 1,32
$row19[0,0] = 77.650000000000006
$row19[0,1] = 66.790000000000006
$row19[0,2] = 55.02
$row19[0,3] = 69.28
$row19[0,4] = 83.93
$row19[0,5] = 97.23
$row19[0,6] = 89.56
$row19[0,7] = 61.49
$row19[0,8] = 30.83
$row19[0,9] = 48.52
$row19[0,10] = 61.75
$row19[0,11] = 29.21
$row19[0,12] = 27.35
$row19[0,13] = 21.6
$row19[0,14] = 53.7
$row19[0,15] = 35.9
$row19[0,16] = 68.599999999999994
$row19[0,17] = 49.28
$row19[0,18] = 41.07
$row19[0,19] = 41.8
$row19[0,20] = 74.61
$row19[0,21] = 43.25
$row19[0,22] = 84.54
$row19[0,23] = 46.54
$row19[0,24] = 42.71
$row19[0,25] = 102.9
$row19[0,26] = 70.099999999999994
$row19[0,27] = 87.38
$row19[0,28] = 68.239999999999995
$row19[0,29] = 50.61
$row19[0,30] = 82.3
$row19[0,31] = 92.51
$ws.Range("B19:AG19").Value = $row19

# Row 20: PARANÁ CLUB DE MUNIQUE
$ws.Range("A20").Value = "PARANÁ CLUB DE MUNIQUE"
$row20 = New-Object 'object[,]' 1,32
$row20[0,0] = 60.59
$row20[0,1] = 82.62
$row20[0,2] = 82.77
$row20[0,3] = 67.02
$row20[0,4] = 56.02
$row20[0,5] = 55.05
$row20[0,6] = 92.96
$row20[0,7] = 90
$row20[0,8] = 40.369999999999997
$row20[0,9] = 16.149999999999999
$row20[0,10] = 87.25
$row20[0,11] = 38.31
$row20[0,12] = 61.3
$row20[0,13] = 31.7
$row20[0,14] = 49.9
$row20[0,15] = 20.399999999999999
$row20[0,16] = 60.3
$row20[0,17] = 43.2
$row20[0,18] = 33.21
$row20[0,19] = 23.7
$row20[0,20] = 63.35
$row20[0,21] = 40.770000000000003
$row20[0,22] = 23.69
$row20[0,23] = 41.54
$row20[0,24] = 30.01
$row20[0,25] = 20.55
$row20[0,26] = 38.56
$row20[0,27] = 26.55
$row20[0,28] = 25.01
$row20[0,29] = 57.27
$row20[0,30] = 65.05
$row20[0,31] = 52.86
$ws.Range("B20:AG20").Value = $row20

# Row 21: CRICIUMA DORTMUND JGS
$ws.Range("A21").Value = "CRICIUMA DORTMUND JGS "
$row21 = New-Object 'object[,]' 1,32
$row21[0,0] = 55.43
$row21[0,1] = 65.95
$row21[0,2] = 82.15
$row21[0,3] = 72.08
$row21[0,4] = 70.930000000000007
$row21[0,5] = 65.06
$row21[0,6] = 74.69
$row21[0,7] = 57.09
$row21[0,8] = 66.27
$row21[0,9] = 55.27
$row21[0,10] = 65.010000000000005
$row21[0,11] = 61.18
$row21[0,12] = 25.07
$row21[0,13] = 39.450000000000003
$row21[0,14] = 36.799999999999997
$row21[0,15] = 23.4
$row21[0,16] = 22.51
$row21[0,17] = 45.79
$row21[0,18] = 48.58
$row21[0,19] = 56.46
$row21[0,20] = 27.86
$row21[0,21] = 42
$row21[0,22] = 31.14
$row21[0,23] = 40.6
$row21[0,24] = 35.74
$row21[0,25] = 28.98
$row21[0,26] = 26.09
$row21[0,27] = 31.44
$row21[0,28] = 14.09
$row21[0,29] = 40.159999999999997
$row21[0,30] = 35.35
$row21[0,31] = 39.18
$ws.Range("B21:AG21").Value = $row21

# AE2 loses its former one-off fill-border style (cellXfs index 6) in favour of the shared border-only style (index 4).
$ws.Range("AE2").Interior.Pattern = -4142

# Restore the active selection to the cell the author left selected.
$ws.Range("U29").Select()
